$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the K column (column G) values for rows 2-7
# (was previously "Strike#"-based values, now recalculated K values)
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 9
